$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value2 = 3249.6667
$ws.Range("J43").Value2 = 2750
$ws.Range("L43").Value2 = 2750
$ws.Range("N43").Value2 = -2888

$ws.Range("H62").Value2 = 7284.4287
$ws.Range("I62").Value2 = 6498
$ws.Range("K62").Value2 = 6498
$ws.Range("M62").Value2 = -5874

$ws.Range("H65").Value2 = 7284.4287
$ws.Range("I65").Value2 = 6498
$ws.Range("K65").Value2 = 32490
$ws.Range("M65").Value2 = -29370

$ws.Range("H86").Value2 = 1233.3334
$ws.Range("I86").Value2 = 350
$ws.Range("K86").Value2 = 350
$ws.Range("M86").Value2 = 773

$ws.Range("H89").Value2 = 1233.3334
$ws.Range("I89").Value2 = 350
$ws.Range("K89").Value2 = 1750
$ws.Range("M89").Value2 = 3866

$ws.Range("H103").Value2 = 2000
$ws.Range("J103").Value2 = 2000
$ws.Range("L103").Value2 = 6000
$ws.Range("N103").Value2 = -7172

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 0
$ws.Range("I2").Value2 = 0
$ws.Range("K2").Value2 = 0
$ws.Range("M2").ClearContents()

$ws.Range("H28").Value2 = 7500
$ws.Range("I28").Value2 = 7500
$ws.Range("K28").Value2 = 7500
$ws.Range("M28").Value2 = -7308

$ws.Range("H32").Value2 = 1628.6765
$ws.Range("I32").Value2 = 1375
$ws.Range("K32").Value2 = 1375
$ws.Range("M32").Value2 = -1088

$ws.Range("H74").Value2 = 2004.8
$ws.Range("I74").Value2 = 2006
$ws.Range("K74").Value2 = 2006
$ws.Range("M74").Value2 = -1132

$ws.Range("H77").Value2 = 2004.8
$ws.Range("I77").Value2 = 2006
$ws.Range("K77").Value2 = 10030
$ws.Range("M77").Value2 = -5662

$ws.Range("H98").Value2 = 0
$ws.Range("J98").Value2 = 0
$ws.Range("L98").Value2 = 0
$ws.Range("N98").ClearContents()

$ws.Range("H99").Value2 = 7500
$ws.Range("I99").Value2 = 7500
$ws.Range("K99").Value2 = 7500
$ws.Range("M99").Value2 = -4505

$ws.Range("H102").Value2 = 2900
$ws.Range("I102").Value2 = 2900
$ws.Range("J102").Value2 = 0
$ws.Range("K102").Value2 = 2900
$ws.Range("L102").Value2 = 0
$ws.Range("M102").Value2 = -1278
$ws.Range("N102").ClearContents()

$ws.Range("H106").Value2 = 28000
$ws.Range("J106").Value2 = 28000
$ws.Range("L106").Value2 = 28000
$ws.Range("N106").Value2 = -30524

$ws.Range("H107").Value2 = 281666.66
$ws.Range("J107").Value2 = 281666.66
$ws.Range("L107").Value2 = 281666.66
$ws.Range("N107").Value2 = -289346.66

$ws.Range("H116").Value2 = 0
$ws.Range("I116").Value2 = 0
$ws.Range("K116").Value2 = 0
$ws.Range("M116").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 0
$ws.Range("I3").Value2 = 0
$ws.Range("K3").Value2 = 0
$ws.Range("M3").ClearContents()

$ws.Range("H86").Value2 = 2745.8096
$ws.Range("I86").Value2 = 2840.1428
$ws.Range("J86").Value2 = 2557.1428
$ws.Range("K86").Value2 = 2840.1428
$ws.Range("L86").Value2 = 2557.1428
$ws.Range("M86").Value2 = -1717.1428
$ws.Range("N86").Value2 = -4803.1428

$ws.Range("H89").Value2 = 2745.8096
$ws.Range("I89").Value2 = 2840.1428
$ws.Range("J89").Value2 = 2557.1428
$ws.Range("K89").Value2 = 14200.714
$ws.Range("L89").Value2 = 12785.714
$ws.Range("M89").Value2 = -8584.714
$ws.Range("N89").Value2 = -24017.714

$ws.Range("H99").Value2 = 2879.4211
$ws.Range("I99").Value2 = 2753.5293
$ws.Range("J99").Value2 = 3949.5
$ws.Range("K99").Value2 = 2753.5293
$ws.Range("L99").Value2 = 3949.5
$ws.Range("M99").Value2 = -1255.5293
$ws.Range("N99").Value2 = -6945.5

$ws.Range("H105").Value2 = 4043.5
$ws.Range("I105").Value2 = 4043.5
$ws.Range("K105").Value2 = 4043.5
$ws.Range("M105").Value2 = -2296.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 3000
$ws.Range("I31").Value2 = 3000
$ws.Range("K31").Value2 = 3000
$ws.Range("M31").Value2 = -2705

$ws.Range("H34").Value2 = 3000
$ws.Range("I34").Value2 = 3000
$ws.Range("K34").Value2 = 3000
$ws.Range("M34").Value2 = -2798

$ws.Range("H58").Value2 = 3508.2727
$ws.Range("I58").Value2 = 1849.4
$ws.Range("J58").Value2 = 4890.6665
$ws.Range("K58").Value2 = 1849.4
$ws.Range("L58").Value2 = 4890.6665
$ws.Range("M58").Value2 = -1646.4
$ws.Range("N58").Value2 = -5296.6665

$ws.Range("H105").Value2 = 958.6875
$ws.Range("I105").Value2 = 874
$ws.Range("J105").Value2 = 1099.8334
$ws.Range("K105").Value2 = 874
$ws.Range("L105").Value2 = 1099.8334
$ws.Range("M105").Value2 = 873
$ws.Range("N105").Value2 = -4593.8334

$ws.Range("H122").Value2 = 1181.4
$ws.Range("J122").Value2 = 961
$ws.Range("L122").Value2 = 2883
$ws.Range("N122").Value2 = -7783

$ws.Range("H134").Value2 = 1211.9
$ws.Range("I134").Value2 = 1088.625
$ws.Range("J134").Value2 = 1705
$ws.Range("K134").Value2 = 3265.875
$ws.Range("L134").Value2 = 5115
$ws.Range("M134").Value2 = -730.875
$ws.Range("N134").Value2 = -10185

$ws.Range("H136").Value2 = 3508.2727
$ws.Range("I136").Value2 = 1849.4
$ws.Range("J136").Value2 = 4890.6665
$ws.Range("K136").Value2 = 5548.200000000001
$ws.Range("L136").Value2 = 14671.9995
$ws.Range("M136").Value2 = -2998.200000000001
$ws.Range("N136").Value2 = -19771.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value2 = 191.66667
$ws.Range("J9").Value2 = 87.5
$ws.Range("L9").Value2 = 262.5
$ws.Range("N9").Value2 = -710.5

$ws.Range("H39").Value2 = 4306.6924
$ws.Range("J39").Value2 = 4306.6924
$ws.Range("L39").Value2 = 12920.0772
$ws.Range("N39").Value2 = -13508.0772

$ws.Range("H128").Value2 = 324197.75
$ws.Range("I128").Value2 = 324197.75
$ws.Range("K128").Value2 = 972593.25
$ws.Range("M128").Value2 = -967613.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value2 = 1252
$ws.Range("I5").Value2 = 504
$ws.Range("J5").Value2 = 2000
$ws.Range("K5").Value2 = 504
$ws.Range("L5").Value2 = 2000
$ws.Range("M5").Value2 = -392
$ws.Range("N5").Value2 = -2224

$ws.Range("H126").Value2 = 16768.5
$ws.Range("I126").Value2 = 11870.333
$ws.Range("J126").Value2 = 21666.666
$ws.Range("K126").Value2 = 35610.999
$ws.Range("L126").Value2 = 64999.99800000001
$ws.Range("M126").Value2 = -33140.999
$ws.Range("N126").Value2 = -69939.99800000001

$ws.Range("H132").Value2 = 2007.1666
$ws.Range("I132").Value2 = 1808.9
$ws.Range("J132").Value2 = 2998.5
$ws.Range("K132").Value2 = 5426.700000000001
$ws.Range("L132").Value2 = 8995.5
$ws.Range("M132").Value2 = -2896.700000000001
$ws.Range("N132").Value2 = -14055.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value2 = 1977.6
$ws.Range("I61").Value2 = 2172
$ws.Range("J61").Value2 = 1200
$ws.Range("K61").Value2 = 2172
$ws.Range("L61").Value2 = 1200
$ws.Range("M61").Value2 = -1970
$ws.Range("N61").Value2 = -1604

$ws.Range("H107").Value2 = 7595
$ws.Range("I107").Value2 = 7595
$ws.Range("K107").Value2 = 7595
$ws.Range("M107").Value2 = -5675

$ws.Range("H113").Value2 = 1977.6
$ws.Range("I113").Value2 = 2172
$ws.Range("J113").Value2 = 1200
$ws.Range("K113").Value2 = 2172
$ws.Range("L113").Value2 = 1200
$ws.Range("M113").Value2 = -2
$ws.Range("N113").Value2 = -5540

$ws.Range("H136").Value2 = 6999.6665
$ws.Range("I136").Value2 = 5999.5
$ws.Range("K136").Value2 = 17998.5
$ws.Range("M136").Value2 = -15448.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value2 = 2475
$ws.Range("I62").Value2 = 2450
$ws.Range("K62").Value2 = 2450
$ws.Range("M62").Value2 = -1826

$ws.Range("H65").Value2 = 2475
$ws.Range("I65").Value2 = 2450
$ws.Range("K65").Value2 = 12250
$ws.Range("M65").Value2 = -9130

$ws.Range("H81").Value2 = 1024.75
$ws.Range("J81").Value2 = 999.5
$ws.Range("L81").Value2 = 1999
$ws.Range("N81").Value2 = -4121

$ws.Range("H84").Value2 = 1024.75
$ws.Range("J84").Value2 = 999.5
$ws.Range("L84").Value2 = 9995
$ws.Range("N84").Value2 = -20603

$ws.Range("H132").Value2 = 2750
$ws.Range("I132").Value2 = 2750
$ws.Range("J132").Value2 = 0
$ws.Range("K132").Value2 = 8250
$ws.Range("L132").Value2 = 0
$ws.Range("M132").Value2 = -5720
$ws.Range("N132").ClearContents()

Write-Output "done"